# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook tracks, per employee, one row per overdue payroll-contribution
# period ("Periodo Mora") together with the amount owed for that period
# ("Valor Mora"). This edit:
#   1) adds one more period (2508) to the table,
#   2) re-sorts the period rows into chronological (ascending) order, and
#   3) makes sure every row's "Valor Mora" matches the rate that actually
#      applied for that period (<=1808 -> 24640, >=1809 -> 31249).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row for the extra period (2508) right after the last
#    existing data row (124), pushing the signature-block rows down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(125).Insert()

# The "last row of the table" formatting used to live on row 124; move it to
# the new last row (125), and restore the regular data-row formatting
# (copied from row 123) onto row 124.
$ws.Range("B124:J124").Copy()
$ws.Range("B125:J125").PasteSpecial(-4122)
$ws.Range("B123:J123").Copy()
$ws.Range("B124:J124").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Build the full chronological period list: 1607..1612, 1701..2412,
#    2501..2508 (110 periods in total -> rows 16..125).
# ---------------------------------------------------------------------------
$periods = @()
foreach ($m in 7..12) { $mm = "{0:D2}" -f $m; $periods += "16$mm" }
foreach ($y in 17..24) {
    $yy = "{0:D2}" -f $y
    foreach ($m in 1..12) { $mm = "{0:D2}" -f $m; $periods += "$yy$mm" }
}
foreach ($m in 1..8) { $mm = "{0:D2}" -f $m; $periods += "25$mm" }

# ---------------------------------------------------------------------------
# 3. Write each row's period + the rate that actually applied to it.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $periods.Count; $i++) {
    $r = 16 + $i
    $p = $periods[$i]
    if ($p -le "1808") { $rate = 24640 } else { $rate = 31249 }
    $ws.Range("E$r").Value = $p
    $ws.Range("F$r").Value = $rate
}

# The newly inserted row also needs its worker/salary columns filled in
# (identical to every other row for this employee).
$ws.Range("B125").Value = "CC"
$ws.Range("C125").Value = "45758329"
$ws.Range("D125").Value = "MARIA VICTORIA NUÑEZ HERNANDEZ"
$ws.Range("G125").Value = 781242

# ---------------------------------------------------------------------------
# 4. Refresh the summary header cells: total "Valor Mora" (E11) and the
#    period count (F13).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 3265556
$ws.Range("F13").Value = 110
